$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Candidatures")

$ws.Range("A5").Value = "Denise D"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = "PHY2710, PHY2710, PHY2400, PHY2400"
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = "Plasmas"
$ws.Range("G5").Value = 3.42
$ws.Range("H5").Value = "PHY2710, PHY2400"

$ws.Activate()
$ws.Range("A6").Select()
